$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Column D/E updates
$ws.Range("D2").Value = '66.710.28'
$ws.Range("E2").Value = '  +4.65%  '
$ws.Range("D3").Value = '3.492.49'
$ws.Range("E3").Value = '  +3.03%  '
Set-TextValue "D5" '592.35'
$ws.Range("E5").Value = '  +3.99%  '
Set-TextValue "D6" '169.39'
$ws.Range("E6").Value = '  +4.78%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.491.60'
$ws.Range("E8").Value = '  +3.01%  '
Set-TextValue "D9" '0.595'
$ws.Range("E9").Value = '  +9.31%  '
$ws.Range("E10").Value = '  +0.48%  '
Set-TextValue "D11" '0.128'
$ws.Range("E11").Value = '  +7.77%  '
$ws.Range("E12").Value = '  +4.71%  '
$ws.Range("D13").Value = '4.094.92'
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("E14").Value = '  -0.15%  '
Set-TextValue "D15" '28.18'
$ws.Range("E15").Value = '  +5.01%  '
$ws.Range("E16").Value = '  +4.43%  '
$ws.Range("D17").Value = '66.704.96'
$ws.Range("E17").Value = '  +4.60%  '
$ws.Range("D18").Value = '3.491.56'
$ws.Range("E18").Value = '  +2.61%  '
$ws.Range("E19").Value = '  +3.84%  '
Set-TextValue "D20" '14.07'
$ws.Range("E20").Value = '  +4.32%  '
Set-TextValue "D21" '392.99'
$ws.Range("E21").Value = '  +4.62%  '
Set-TextValue "D22" '7.93'
$ws.Range("E22").Value = '  +2.39%  '
Set-TextValue "D23" '73.12'
$ws.Range("E23").Value = '  +4.43%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  +5.10%  '
$ws.Range("E26").Value = '  +7.58%  '
Set-TextValue "D27" '10.25'
$ws.Range("E27").Value = '  +8.01%  '
$ws.Range("E28").Value = '  +1.81%  '
Set-TextValue "D29" '1.00'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  +4.80%  '
$ws.Range("E31").Value = '  +6.30%  '
$ws.Range("E32").Value = '  +3.72%  '
Set-TextValue "D33" '23.55'
$ws.Range("E34").Value = '  +6.14%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  +10.33%  '
Set-TextValue "D37" '161.67'
$ws.Range("E37").Value = '  +1.33%  '
Set-TextValue "D38" '0.902'
$ws.Range("E38").Value = '  +5.29%  '
$ws.Range("E39").Value = '  +7.13%  '
Set-TextValue "D44" '26.86'
$ws.Range("E44").Value = '  +3.24%  '
Set-TextValue "D45" '43.23'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("D46").Value = '2.767.07'
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("E48").Value = '  +4.23%  '
Set-TextValue "D49" '346.91'
$ws.Range("E49").Value = '  +5.75%  '
$ws.Range("E50").Value = '  +5.96%  '
Set-TextValue "D51" '0.889'
$ws.Range("E51").Value = '  +9.72%  '

# Rows 40-43: coin order swap (RenderToken/Hedera moved up, Filecoin/EnergySwap moved down)
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D40" '6.77'
$ws.Range("E40").Value = '  +6.11%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D41" '0.0745'
$ws.Range("E41").Value = '  +3.70%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D42" '4.66'
$ws.Range("E42").Value = '  +7.43%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D43" '26.66'
$ws.Range("E43").Value = '  +4.27%  '
